$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gmail_Signup_Testdata")

# Replace the test data row (row 2): FirstName/LastName/Password values.
# The leading apostrophe on C2 preserves the existing "quote prefix" text
# formatting that was already applied to that cell.
$ws.Range("C2").Value = "'Mark"
$ws.Range("D2").Value = "Peterson"
$ws.Range("E2").Value = "January@2021#"

# Leave the same cell selected as in the authored workbook.
$ws.Range("E7").Select()
